# FORM.docx update — "feat: supposed to be finish"
#
# 1. The department/unit placeholder cell changes from a literal value to
#    a template tag.
# 2/3. Two signatory names that used to be split across several runs
#    (because Word's spell-checker had flagged individual words, wrapping
#    them in <w:proofErr> start/end markers) are retyped as a single run
#    of plain text, which naturally drops the now-stale proofErr markers.
# 4. The {nama_pemohon} placeholder loses its spell-check proofErr
#    wrapper around "nama_pemohon", while keeping "{", "nama_pemohon" and
#    "}" as three separate runs.

$d = $word.ActiveDocument

# 1. "KP & TI" -> "{nama_bagian}" (department/unit table cell)
$d.Content.Find.Execute("KP & TI", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{nama_bagian}", 2)

# 2. Re-type "Andi Nurhasbi Alauddin, S.E., M.H." as one clean run —
#    text is unchanged, but this merges the runs that spell-check had
#    split and clears the associated proofErr markers.
$d.Content.Find.Execute("Andi Nurhasbi Alauddin, S.E., M.H.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Andi Nurhasbi Alauddin, S.E., M.H.", 2)

# 3. Same treatment for "Muhammad Andy Alfariz, A.Md.Ak."
$d.Content.Find.Execute("Muhammad Andy Alfariz, A.Md.Ak.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Muhammad Andy Alfariz, A.Md.Ak.", 2)

# 4. Clear the stale proofErr spell-check wrapper around the
#    {nama_pemohon} placeholder. First re-type the whole tag as a single
#    run (this is what actually drops the <w:proofErr> markers), then
#    restore the original "{" / "nama_pemohon" / "}" run split by
#    touching (and immediately releasing) character formatting on the
#    inner span only — this forces Word to split the run again at the
#    same two boundaries without altering the visible text or format.
$d.Content.Find.Execute("{nama_pemohon}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{nama_pemohon}", 2)

$rng = $d.Content
$rng.Find.Execute("{nama_pemohon}")
$start = $rng.Start
$end = $rng.End
$inner = $d.Range($start + 1, $end - 1)
$inner.Font.Bold = $true

$inner2 = $d.Content
$inner2.Find.Execute("nama_pemohon")
$inner2.Font.Bold = $false
